$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Correct the vaccination count for 2021-05-31 (row 120, date 44305):
#    3001 -> 3672. Everything else in that row is formula-driven and
#    will recompute automatically.
# ---------------------------------------------------------------------
$ws.Range("B120").Value = 3672
$ws.Rows(120).RowHeight = 13

# ---------------------------------------------------------------------
# 2. Make room for a new day of data (2021-06-02, date 44307) by
#    inserting a fresh row right above the current last row (121).
#    Copying row 120 down and inserting it pushes the old row 121 to
#    122 (carrying its formatting + formulas with it, relative refs
#    auto-adjusted by Excel) and seeds row 121 with row 120's own
#    formatting/formula pattern.
# ---------------------------------------------------------------------
$ws.Rows(120).Copy()
$ws.Rows(121).Insert()

# Row 121 now holds a shifted copy of row 120's data; overwrite with
# the real values for 2021-06-01 (date 44306), matching what used to
# be the last row before this edit.
$ws.Range("A121").Value = 44306
$ws.Range("B121").Value = 3058
$ws.Range("H121").Value = 28612
$ws.Range("J121").Value = 7153

# ---------------------------------------------------------------------
# 3. Row 122 (the old last row, shifted down intact) now needs to
#    become the real new last row: 2021-06-02 (date 44307).
# ---------------------------------------------------------------------
$ws.Range("A122").Value = 44307
$ws.Range("B122").Value = 3698
$ws.Range("J122").Value = 7153
$ws.Range("C122").Formula = "=(AVERAGE(B116:B122))"
$ws.Range("D122").Formula = "=(D121-B122)"
$ws.Range("E122").Formula = "=E121+B122"
$ws.Range("I122").Formula = "=G122/2"

# ---------------------------------------------------------------------
# 4. Re-enter the F/K/L formulas across the whole 116-122 block as one
#    continuous fill so Excel consolidates them into shared-formula
#    groups the same way a drag-fill / copy-down would.
# ---------------------------------------------------------------------
$ws.Range("F116:F122").Formula = "=(E116-G116)"
$ws.Range("K116:K121").Formula = "=D116/C116"
$ws.Range("L116:L121").Formula = "=A116+K116"

# Row 122's own K/L formulas (kept as their own single-cell pattern,
# matching how every other "latest row" has historically been added).
$ws.Range("K122").Formula = "=D122/C122"
$ws.Range("L122").Formula = "=A122+K122"

# ---------------------------------------------------------------------
# 5. Selection / view bookkeeping to match where the author left off.
# ---------------------------------------------------------------------
$ws.Range("H132").Select()

$excel.CalculateFull()
